# Q3 Update - 2025
# Applies the following changes to the "fromCSV" worksheet:
#  1. Updates the asylum_seekers value (column O) for row 60 (Cuba/2024)
#     from 8 to 7.
#  2. Updates the asylum_seekers value (column O) for row 61 (Jamaica/2024)
#     from 9 to 8.
#  3. Updates the short-url value (column B) for all data rows from
#     "xeU9S4" to "sW71rO".
#  4. Removes the last data row (row 63 - Venezuela/2024) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row so the short-url replacement covers every
# data row regardless of how many rows currently exist.
$lastRow = $ws.UsedRange.Rows.Count

# 1. & 2. Update the asylum_seekers figures for the Cuba and Jamaica rows.
#    These columns store every value as text (the source data is a CSV
#    export where even numeric-looking cells are shared strings), so the
#    new figures are copied in from other cells in the same column that
#    already hold the desired text ("7" / "8") instead of being typed in
#    directly - a plain numeric assignment would be auto-converted to a
#    real number by Excel and change the cell's underlying type/style.
$ws.Range("O34").Copy($ws.Range("O60"))
$ws.Range("O11").Copy($ws.Range("O61"))

# 3. Replace the short-url shared value in column B for all data rows
#    (row 1 is the header and must stay untouched).
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq "xeU9S4") {
        $cell.Value = "sW71rO"
    }
}

# 4. Remove the Venezuela row (last row in the sheet) entirely, shifting
#    nothing up beneath it since it is the final row of data.
$ws.Rows(63).Delete()
